$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "abbr" header in column E
$ws.Range("E1").Value = "abbr"

# Fill in abbreviation values for the epoch rows that have a color (D column) set:
# row 13 Late Oligocene -> Late
# row 14 Early Oligocene -> Early
# row 15 Late Eocene -> Late
# row 16 Middle Eocene -> Mid
# row 17 Early Eocene -> Early
$ws.Range("E13").Value = "Late"
$ws.Range("E14").Value = "Early"
$ws.Range("E15").Value = "Late"
$ws.Range("E16").Value = "Mid"
$ws.Range("E17").Value = "Early"
